# ListaVecinos.xlsx update
# - House "25" now also covers unit 89 -> "25/89"
# - Neighbor at house 67 changed from "Lopez Omar" to "Lopez Guiliana"
# - House range "72 / 73" extended to include 88 -> "72 / 73 / 88"
# - Neighbor "Maria " (with trailing space, surname missing) completed to "Maria Tapia"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A17").Value = "25/89"
$ws.Range("B53").Value = "Lopez Guiliana"
$ws.Range("A58").Value = "72 / 73 / 88"
$ws.Range("B68").Value = "Maria Tapia"

# Match the workbook's updated view/selection state
$ws.Activate()
$ws.Range("A59").Select()
